$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 173; this shifts the existing
# rows 173-208 down to 174-209 (carrying their formatting/styles with them),
# and the workbook's used-range dimension grows from A1:R208 to A1:R209.
$ws.Rows("173:173").Insert()

# Populate the newly-inserted row 173 with the new weekly price entry.
$ws.Range("A173").Value = 9
$ws.Range("B173").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C173").Value = "Metropolitana"
$ws.Range("D173").Value = 44641
$ws.Range("E173").Value = 13
$ws.Range("F173").Value = 100112003
$ws.Range("G173").Value = "Ajo"
$ws.Range("H173").Value = "Chino"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 520
$ws.Range("K173").Value = 17500
$ws.Range("L173").Value = 18000
$ws.Range("M173").Value = 17750
$ws.Range("N173").Value = "`$/caja 10 kilos"
$ws.Range("O173").Value = "China"
$ws.Range("P173").Value = 1775
$ws.Range("Q173").Value = 10
$ws.Range("R173").Value = "Hortaliza"
